$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocalData")

$ws.Range("A49").Value = "MapNpcData.MapNpcMenu.1000023.1"
$ws.Range("B49").Value = "마이홈 돌아가기"

# Touch C49/D49 so they persist as explicit (empty) cells in the sheet,
# matching the empty inlineStr cells present in every other data row.
$ws.Range("C49").Font.Bold = $false
$ws.Range("D49").Font.Bold = $false
